$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Förändrad" column (C) records the date the row was last refreshed,
# stored as an Excel date serial number. The automatic update run bumped
# this value by one day (45171 -> 45172 = 2023-09-02 -> 2023-09-03) for
# every data row (rows 2 through 439).
$ws.Range("C2:C439").Value = 45172
